# Auto-generated Excel COM-interop script implementing the weekly CompStat data refresh
# (report number 42 -> 43, week range 10/14-10/20/2024 -> 10/21-10/27/2024,
# and the updated crime-statistics figures for rows 14-31, columns C:N).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header text: report Volume/Number and the week-covering date range ---
# These strings live inside rich-text shared strings; updating the whole cell
# text (each is the concatenation of several runs) reproduces the same visible text
# while preserving formatting, since the runs share identical rPr in this cell.
$ws.Range("C8").Value = "Volume 31   Number  43"
$ws.Range("C9").Value = "Report Covering the Week  10/21/2024  Through  10/27/2024"

# --- 2. Fix up cells whose underlying type switches between text placeholder ---
#        ("0" / "***.*") and a real number, so the number format / style matches
#        the surrounding data cells. We copy from a donor cell that already has
#        the desired style+content, then (for text->number) overwrite the value.

# -- text placeholder -> real number (style becomes 14 count / 15 percent) --
$ws.Range("D16").Copy($ws.Range("D15"))
$ws.Range("E16").Copy($ws.Range("E15"))
$ws.Range("C15").Copy($ws.Range("C23"))
$ws.Range("D16").Copy($ws.Range("D23"))
$ws.Range("E16").Copy($ws.Range("E23"))
$ws.Range("F15").Copy($ws.Range("F23"))
$ws.Range("D16").Copy($ws.Range("D27"))
$ws.Range("E16").Copy($ws.Range("E27"))
$ws.Range("C15").Copy($ws.Range("C31"))

# -- real number -> text placeholder (style becomes 13, General) --
$ws.Range("D14").Copy($ws.Range("D31"))
$ws.Range("E14").Copy($ws.Range("E31"))

# --- 3. Assign the new numeric values (weekly/28-day/YTD counts and %-changes) ---
# Row 14
$ws.Range("M14").Value = -33.333333333333
# Row 15
$ws.Range("C15").Value = 2
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 100
$ws.Range("F15").Value = 3
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 17
$ws.Range("J15").Value = 11
$ws.Range("K15").Value = 54.545454545454
$ws.Range("L15").Value = 6.25
$ws.Range("M15").Value = 70
$ws.Range("N15").Value = 6.25
# Row 16
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 50
$ws.Range("F16").Value = 11
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = -31.25
$ws.Range("I16").Value = 130
$ws.Range("J16").Value = 149
$ws.Range("K16").Value = -12.751677852349
$ws.Range("L16").Value = -20.245398773006
$ws.Range("M16").Value = -10.958904109589
$ws.Range("N16").Value = -84.848484848484
# Row 17
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 66.666666666666
$ws.Range("F17").Value = 14
$ws.Range("G17").Value = 21
$ws.Range("H17").Value = -33.333333333333
$ws.Range("I17").Value = 194
$ws.Range("J17").Value = 193
$ws.Range("K17").Value = 0.518134715025
$ws.Range("L17").Value = -3.960396039603
$ws.Range("M17").Value = 67.241379310344
$ws.Range("N17").Value = -18.143459915611
# Row 18
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 50
$ws.Range("I18").Value = 165
$ws.Range("J18").Value = 185
$ws.Range("K18").Value = -10.810810810810
$ws.Range("L18").Value = -31.535269709543
$ws.Range("M18").Value = -20.673076923076
$ws.Range("N18").Value = -87.631184407796
# Row 19
$ws.Range("C19").Value = 19
$ws.Range("D19").Value = 19
$ws.Range("E19").Value = 0
$ws.Range("G19").Value = 87
$ws.Range("H19").Value = -14.942528735632
$ws.Range("I19").Value = 856
$ws.Range("J19").Value = 861
$ws.Range("K19").Value = -0.580720092915
$ws.Range("L19").Value = -0.116686114352
$ws.Range("M19").Value = -28.308207705192
$ws.Range("N19").Value = -62.039911308204
# Row 20
$ws.Range("F20").Value = 3
$ws.Range("H20").Value = -40
$ws.Range("I20").Value = 36
$ws.Range("J20").Value = 54
$ws.Range("K20").Value = -33.333333333333
$ws.Range("L20").Value = -43.75
$ws.Range("M20").Value = -10
$ws.Range("N20").Value = -96.491228070175
# Row 21
$ws.Range("C21").Value = 36
$ws.Range("D21").Value = 30
$ws.Range("E21").Value = 20
$ws.Range("F21").Value = 123
$ws.Range("G21").Value = 148
$ws.Range("H21").Value = -16.891891891891
$ws.Range("I21").Value = 1400
$ws.Range("J21").Value = 1454
$ws.Range("K21").Value = -3.713892709766
$ws.Range("L21").Value = -9.385113268608
$ws.Range("M21").Value = -18.462434478742
$ws.Range("N21").Value = -75.575715282623
# Row 22
$ws.Range("C22").Value = 4
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 300
$ws.Range("F22").Value = 8
$ws.Range("G22").Value = 11
$ws.Range("H22").Value = -27.272727272727
$ws.Range("I22").Value = 54
$ws.Range("J22").Value = 77
$ws.Range("K22").Value = -29.870129870129
$ws.Range("L22").Value = -28.947368421052
$ws.Range("M22").Value = -3.571428571428
# Row 23
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 1
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 8
$ws.Range("J23").Value = 10
$ws.Range("K23").Value = -20
$ws.Range("L23").Value = -46.666666666666
$ws.Range("M23").Value = -46.666666666666
# Row 24
$ws.Range("C24").Value = 57
$ws.Range("D24").Value = 44
$ws.Range("E24").Value = 29.545454545454
$ws.Range("G24").Value = 167
$ws.Range("H24").Value = 41.916167664670
$ws.Range("I24").Value = 2445
$ws.Range("J24").Value = 1810
$ws.Range("K24").Value = 35.082872928176
$ws.Range("L24").Value = 24.490835030549
$ws.Range("M24").Value = 59.908436886854
# Row 25
$ws.Range("C25").Value = 46
$ws.Range("D25").Value = 37
$ws.Range("E25").Value = 24.324324324324
$ws.Range("F25").Value = 186
$ws.Range("G25").Value = 124
$ws.Range("H25").Value = 50
$ws.Range("I25").Value = 2073
$ws.Range("J25").Value = 1344
$ws.Range("K25").Value = 54.241071428571
$ws.Range("L25").Value = 38.384512683578
# Row 26
$ws.Range("C26").Value = 15
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = 114.285714285714
$ws.Range("F26").Value = 50
$ws.Range("G26").Value = 27
$ws.Range("H26").Value = 85.185185185185
$ws.Range("I26").Value = 502
$ws.Range("J26").Value = 420
$ws.Range("K26").Value = 19.523809523809
$ws.Range("L26").Value = 17.016317016317
$ws.Range("M26").Value = 42.209631728045
# Row 27
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 100
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 28
$ws.Range("J27").Value = 21
$ws.Range("K27").Value = 33.333333333333
$ws.Range("L27").Value = 7.692307692307
# Row 28
$ws.Range("C28").Value = 1
$ws.Range("E28").Value = -75
$ws.Range("F28").Value = 11
$ws.Range("G28").Value = 13
$ws.Range("H28").Value = -15.384615384615
$ws.Range("I28").Value = 110
$ws.Range("J28").Value = 98
$ws.Range("K28").Value = 12.244897959183
$ws.Range("L28").Value = 22.222222222222
# Row 29
$ws.Range("M29").Value = -33.333333333333
# Row 30
$ws.Range("M30").Value = 0
# Row 31
$ws.Range("C31").Value = 2
$ws.Range("F31").Value = 3
$ws.Range("H31").Value = 200
$ws.Range("I31").Value = 9
$ws.Range("K31").Value = 28.571428571428
$ws.Range("L31").Value = -55

Write-Output "Applied weekly CompStat data refresh."
